$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 100 (ID 99): Reverse Nodes in K Group ---
$ws.Range("B99:H99").Copy()
$ws.Range("B100:H100").PasteSpecial($xlPasteFormats)
$ws.Range("B100").Value = "Linked List"
$ws.Range("C100").Value = "Reverse Nodes in K Group"
$ws.Range("D100").Value = "Hard"
$ws.Range("E100").Value = "Done"
$ws.Range("F100").Value = 45925
$ws.Range("G100").Value = "O(n)"
$ws.Range("H100").Value = "O(1)"

# --- Row 101 (ID 100): Rotate List ---
$ws.Range("B99:I99").Copy()
$ws.Range("B101:I101").PasteSpecial($xlPasteFormats)
$ws.Range("B101").Value = "Linked List"
$ws.Range("C101").Value = "Rotate List"
$ws.Range("D101").Value = "Medium"
$ws.Range("E101").Value = "Done"
$ws.Range("F101").Value = 45925
$ws.Range("G101").Value = "O(n)"
$ws.Range("H101").Value = "O(1)"
$ws.Range("I101").Value = "Slow & Fast Pointers"

# --- Row 102 (new row, ID 101): Insertion Sort ---
$ws.Range("A99:I99").Copy()
$ws.Range("A102:I102").PasteSpecial($xlPasteFormats)
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = "Linked List"
$ws.Range("C102").Value = "Insertion Sort"
$ws.Range("D102").Value = "Medium"
$ws.Range("E102").Value = "Done"
$ws.Range("F102").Value = 45925
$ws.Range("G102").Value = "o(n ^ 2)"
$ws.Range("H102").Value = "O(1)"
$ws.Range("I102").Value = "Insertion Sort"

# --- Row 103 (new row, ID 102): Design Linked List ---
$ws.Range("A99:H99").Copy()
$ws.Range("A103:H103").PasteSpecial($xlPasteFormats)
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "Linked List"
$ws.Range("C103").Value = "Design Linked List"
$ws.Range("D103").Value = "Medium"
$ws.Range("E103").Value = "Done"
$ws.Range("F103").Value = 45925
$ws.Range("G103").Value = "O(n)"
$ws.Range("H103").Value = "O(1)"

# --- Rows 104-109 (new, blank placeholder rows, ID only) ---
$ws.Range("A99").Copy()
for ($r = 104; $r -le 109; $r++) {
    $ws.Range("A$r").PasteSpecial($xlPasteFormats)
    $ws.Range("A$r").Value = ($r - 1)
}

# --- Rows 110-111: add the ID column to the existing "To Do" legend rows ---
$ws.Range("A99").Copy()
$ws.Range("A110").PasteSpecial($xlPasteFormats)
$ws.Range("A110").Value = 109

$ws.Range("A99").Copy()
$ws.Range("A111").PasteSpecial($xlPasteFormats)
$ws.Range("A111").Value = 110

# --- Final UI state: selection moved to A107 (matches the author's last edit) ---
$ws.Range("A107").Select()
